$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.347.89"
$ws.Range("E2").Value = "  +12.57%  "
$ws.Range("D3").Value = "1.824.72"
$ws.Range("E3").Value = "  +9.07%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'229.48"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'31.66"
$ws.Range("E8").Value = "  +6.60%  "
$ws.Range("D9").Value = "'47.20"
$ws.Range("E9").Value = "  +6.94%  "
$ws.Range("D10").Value = "'0.284"
$ws.Range("E10").Value = "  +7.43%  "
$ws.Range("D11").Value = "'0.0674"
$ws.Range("E11").Value = "  +5.33%  "
$ws.Range("D12").Value = "'0.0929"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "2.086.21"
$ws.Range("E13").Value = "  +9.01%  "
$ws.Range("D14").Value = "1.837.90"
$ws.Range("E14").Value = "  +9.91%  "
$ws.Range("D15").Value = "'0.648"
$ws.Range("E15").Value = "  +5.34%  "
$ws.Range("D16").Value = "'10.42"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "34.276.30"
$ws.Range("E17").Value = "  +12.27%  "
$ws.Range("D18").Value = "'4.29"
$ws.Range("E18").Value = "  +7.54%  "
$ws.Range("D19").Value = "'69.81"
$ws.Range("E19").Value = "  +5.24%  "
$ws.Range("D20").Value = "'258.77"
$ws.Range("E20").Value = "  +6.52%  "
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  +6.02%  "
$ws.Range("D24").Value = "'4.34"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("D26").Value = "'159.64"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'16.64"
$ws.Range("E27").Value = "  +4.85%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  +7.33%  "
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'3.91"
$ws.Range("E31").Value = "  +12.67%  "
$ws.Range("D32").Value = "'0.0521"
$ws.Range("E32").Value = "  +5.11%  "
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("E34").Value = "  +8.44%  "
$ws.Range("D35").Value = "1.548.77"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").Value = "'1.80"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("E37").Value = "  +6.47%  "
$ws.Range("D38").Value = "'0.633"
$ws.Range("E38").Value = "  +5.51%  "
$ws.Range("E39").Value = "  +6.50%  "
$ws.Range("D40").Value = "'84.76"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").Value = "'0.920"
$ws.Range("E42").Value = "  +9.61%  "
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("D44").Value = "'2.17"
$ws.Range("E44").Value = "  +9.49%  "
$ws.Range("E45").Value = "  +5.93%  "
$ws.Range("D46").Value = "'1.12"
$ws.Range("E46").Value = "  +171.62%  "
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("D48").Value = "1.990.46"
$ws.Range("E48").Value = "  +10.20%  "
$ws.Range("D49").Value = "'12.29"
$ws.Range("E49").Value = "  +27.18%  "
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "  -0.15%  "
